# Add dtau to tables, add long heating graph
#
# - Rename the three existing sheets from "Sheet1" / "Sheet1 (2)" / "Sheet1 (3)"
#   to simply "1" / "2" / "3".
# - Refresh each sheet's selection to F6 (matches the post-edit state).
# - Duplicate the first sheet into a new 4th sheet named "testowy" that models
#   a long heating run: dTau (B6/B7) is dropped from 2E-3/2.5E-3 down to 1E-4/1E-4.
# - Leave the 3rd sheet ("3") as the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- rename the existing sheets -------------------------------------------------
$wb.Worksheets.Item(1).Name = "1"
$wb.Worksheets.Item(2).Name = "2"
$wb.Worksheets.Item(3).Name = "3"

# --- refresh selection on each of the original sheets ---------------------------
foreach ($idx in 1..3) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Activate()
    $ws.Range("F6").Select() | Out-Null
}

# --- add the new "testowy" sheet by duplicating sheet "1" -----------------------
$src = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "testowy"

# dTau values (v0 / v1) for the long-heating run
$newSheet.Range("B6").Value = 0.0001
$newSheet.Range("B7").Value = 0.0001

$newSheet.Activate()
$newSheet.Range("J7").Select() | Out-Null

# --- sheet "3" stays the selected/active tab -------------------------------------
$wb.Worksheets.Item(3).Activate()
